$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the active selection on the sheet view
$ws.Range("J35").Select()

# Rows that need the avgsalary (column C) value inserted, with
# shoparea/foodseats/retailturnover shifted right into D/E/F
# (the previous F value, a stray ratio, is discarded).
$rows = @(
    @{ Row = 17; NewC = 19406.97712 },
    @{ Row = 18; NewC = 17100.190320000002 },
    @{ Row = 19; NewC = 18039.2 },
    @{ Row = 20; NewC = 19598.053680000001 },
    @{ Row = 21; NewC = 24531.17148999999 }
)

foreach ($item in $rows) {
    $r = $item.Row

    $oldC = $ws.Cells.Item($r, 3).Value()
    $oldD = $ws.Cells.Item($r, 4).Value()
    $oldE = $ws.Cells.Item($r, 5).Value()

    $ws.Cells.Item($r, 6).Value = $oldE
    $ws.Cells.Item($r, 5).Value = $oldD
    $ws.Cells.Item($r, 4).Value = $oldC
    $ws.Cells.Item($r, 3).Value = $item.NewC
}
